# Calculadora Solar V4.1 - make string calcs editable + add "Paneles totales reales"
# This script reproduces (in terms of resulting cell values/formulas/formatting)
# the OOXML diff supplied for the commit:
#   "Fix: Hacer strings editables y corregir lógica de paneles totales"

$wb = $excel.ActiveWorkbook
$wsCalc = $wb.Worksheets.Item("Calculadora")
$wsInforme = $wb.Worksheets.Item("Informe Cliente")

# ---------------------------------------------------------------------------
# 1) Rows 40-42: turn the string Voltage/Current/Power formulas into plain,
#    user-editable numbers (matching the new "input" style used elsewhere,
#    e.g. D34/D35/D47/D48/D50), and update their "FÓRMULA" hint text.
# ---------------------------------------------------------------------------

function Set-InputStyle($rng) {
    # Mirrors style s="18" in styles.xml: bold blue text, yellow fill,
    # thin border all round, centered, 2-decimal number format.
    $rng.NumberFormat = "0.00"
    $rng.Font.Bold = $true
    $rng.Font.Size = 11
    $rng.Font.Color = 0x00C07000   # BGR for 0070C0 (bold blue)
    $rng.Interior.Pattern = -4124  # xlSolid
    $rng.Interior.Color = 0x009CF2FF  # BGR for FFF2CC
    $rng.HorizontalAlignment = -4108 # xlCenter
    $rng.Borders.Item(7).LineStyle = 1
    $rng.Borders.Item(10).LineStyle = 1
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(9).LineStyle = 1
}

$d40 = $wsCalc.Range("D40")
$d40.Formula = "73"
Set-InputStyle $d40
$wsCalc.Range("F40").Value = "EDITAR según modelo inversor"

$d41 = $wsCalc.Range("D41")
$d41.Formula = "10.96"
Set-InputStyle $d41
$wsCalc.Range("F41").Value = "EDITAR según modelo inversor"

$d42 = $wsCalc.Range("D42")
$d42.Formula = "800"
Set-InputStyle $d42
$wsCalc.Range("F42").Value = "EDITAR según configuración"

# ---------------------------------------------------------------------------
# 2) Row 60: update the "FÓRMULA" hint text.
# ---------------------------------------------------------------------------
$wsCalc.Range("F60").Value = "Paneles necesarios / Paneles por string"

# ---------------------------------------------------------------------------
# 3) Insert a new row 61 "Paneles totales reales" = D60*D37 (pushes the old
#    rows 61 "Potencia del inversor" -> 62 and 62 "Corriente controlador" ->
#    63 down, which also auto-updates the D58/D27 & MAX(D5:D22) formulas'
#    row numbers automatically via the Insert).
# ---------------------------------------------------------------------------
$wsCalc.Rows.Item(61).Insert()

$wsCalc.Range("B61").Value = "Paneles totales reales"
$wsCalc.Range("D61").Formula = "=D60*D37"
$wsCalc.Range("E61").Value = "unidades"
$wsCalc.Range("F61").Value = "Strings × Paneles/string = Total a comprar"

# Remove the placeholder empty cells the engine created in columns that the
# new row does not use (C/G/H/I/J) so the row only carries B/D/E/F, like the
# target sheet.
$wsCalc.Range("C61").ClearContents()
$wsCalc.Range("G61").ClearContents()
$wsCalc.Range("H61").ClearContents()
$wsCalc.Range("I61").ClearContents()
$wsCalc.Range("J61").ClearContents()

# Style B61 (label, bold left aligned, thin border) / E61 (centered, thin
# border) / F61 (italic small, thin border) like their column neighbours,
# and D61 (green highlight "result" style, but with an integer 0-format).
$b61 = $wsCalc.Range("B61")
$b61.Font.Bold = $true
$b61.HorizontalAlignment = -4131 # xlLeft
$b61.Borders.Item(7).LineStyle = 1
$b61.Borders.Item(10).LineStyle = 1
$b61.Borders.Item(8).LineStyle = 1
$b61.Borders.Item(9).LineStyle = 1

$e61 = $wsCalc.Range("E61")
$e61.HorizontalAlignment = -4108 # xlCenter
$e61.Borders.Item(7).LineStyle = 1
$e61.Borders.Item(10).LineStyle = 1
$e61.Borders.Item(8).LineStyle = 1
$e61.Borders.Item(9).LineStyle = 1

$f61 = $wsCalc.Range("F61")
$f61.Font.Italic = $true
$f61.Font.Size = 9
$f61.Borders.Item(7).LineStyle = 1
$f61.Borders.Item(10).LineStyle = 1
$f61.Borders.Item(8).LineStyle = 1
$f61.Borders.Item(9).LineStyle = 1

$d61 = $wsCalc.Range("D61")
$d61.NumberFormat = "0"
$d61.Font.Bold = $true
$d61.Font.Size = 11
$d61.Interior.Pattern = -4124   # xlSolid
$d61.Interior.Color = 0x00DAEFE2  # BGR for E2EFDA (light green)
$d61.HorizontalAlignment = -4108  # xlCenter
$d61.Borders.Item(7).LineStyle = 1
$d61.Borders.Item(10).LineStyle = 1
$d61.Borders.Item(8).LineStyle = 1
$d61.Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# 4) "Informe Cliente" sheet: fix cross-sheet references.
#    B24/B27 already auto-shifted (D61->D62, D62->D63) by the row insert;
#    B19 needs to be repointed deliberately to the NEW D61 (paneles reales).
# ---------------------------------------------------------------------------
$wsInforme.Range("B19").Formula = '=Calculadora!D61&" paneles de "&TEXT(Calculadora!D33,"0.00")&"W"'

Write-Output "edit applied"
